$d = $word.ActiveDocument

# 1. Update title paragraph: date and paper title
$d.Content.Find.Execute('המאמר היומי של מייק - 17.03.25', $false, $false, $false, $false, $false, $true, 1, $false, 'המאמר היומי של מייק - 15.03.25', 2) | Out-Null
$d.Content.Find.Execute('JanusFlow: Harmonizing Autoregression and Rectified Flow for Unified Multimodal Understanding and Generation', $false, $false, $false, $false, $false, $true, 1, $false, 'Generative Representational Instruction Tuning', 2) | Out-Null

# 2. Replace the body paragraphs (2-6) with the new review text
$d.Content.Find.Execute('מזמן לא סקרתי מאמר על מודלים גנרטיביים מולטימודליים. מודלים אלו מאומנים לא רק לגנרט דאטה מכמה סוגים (במקרה של JanusFlow של שפה טבעית ותמונות) אלא גם לבצע משימות הכרוכות בהבנה של הקשרים ביו המודליות האלו. למשל מודל מולטימודלי בתחום שפה ותמונות צריך להיות לענות על שאלות על תמונה. המודל מורכב ממודל עיקרי (הנקרא LLM) וכמה אנקודרים ודקורדרים המיועדים לייצוג דאטה ממודליות שונות והפיכתו של ייצוגו לפיסת דאטה (דקודרים). כל המודלים במאמר מבוססים על הטרנספורמרים באופן מאוד לא מפתיע.', $false, $false, $false, $false, $false, $true, 1, $false, 'נתקלתי במאמר הזה די במקרה - תוך כדי איזה שיחה עם LLM מצוי על נושא של אמבדינגס הקשריים (contextualized embeddings) ואופן בנייתם. המאמר די קליל וחשבתי שאם כבר השקעתי 5 דקות בקריאתו אז אשקיע עוד 10 דקות בסקירתו. המאמר מציע שיטה המאחדת instruction tuning (נקרא לזה InTn) למטרת גנרוט ו-InTn למטרות בניית ייצוג דאטה הקשרי. ', 2) | Out-Null
$d.Content.Find.Execute('המאמר מציע שיטה לאמן מולטימודלי (הנקרא LLM במאמר) כזה כאשר הפרט המעניין לגביו הוא שימוש באנקודרים שונים לשפה ולתמונות (ברוב המודלים המולטימודליים משתמשים באותו מודל backbone). בגדול במהלך האימון המודל לומד לחזות את הטוקנים של תשובה על פרומפט נתון כאשר פרומפט ותשובה יכולים להיות גם טוקן ויזואלי (ייצוג של פאץ'' של תמונה) וגם טוקן רגיל(= סדרת אותיות). בנוסף הפרומפט יכול להיות שילוב של טוקנים ויזואליים וטוקנים של השפה במשימת visual question answering. בנוסף (לא מופיע במאמר הזה בצורה מפורשת אך נעשה במודלים מולטימודליים אחרים) המודל מאומן גם על דאטה טקסטואלי בלבד(כמו ב-pretraining של מודל שפה רגיל)', $false, $false, $false, $false, $false, $true, 1, $false, 'מטרת InTn גנרטיבי (generative instruction tuning) הוא די מובן ומטרתו לאמן את המודל למלא את הוראות המשתמש (לדוגמא לבניית chatbot). לעומת זאת מטרת InTn ייצוגי (representational instruction tuning) היא לאמן מודל אנקודר, הבונה ייצוג וקטורי של טקסט, בתלות בהוראות המשתמש (שזה די קרוב לייצוג הקשרי). יש לא מעט מאמרים הדנים באיך לפתח מודל המסוגל לבצע כל משימה כזו בנפרד - והמאמר הזה מציע שיטה שמאמנת את אותו המודל לעשות את שני הדברים האלו (לא באותו הזמן אמנם).', 2) | Out-Null
$d.Content.Find.Execute('כמה פרטים על המודלים השונים (פרט ל-LLM) המופיעים במאמר. עבור דאטה שפתי הטוקנים עוברים אנדוקר מאומן (נקרא und enc) -  אחרי הטוקנים עוברים שכבה לינארית מאומנת. עבור דאטה ויזואלי יש אנקודר סטנדרטי לא מאומן המבוסס על VAE ואחרי יש עוד אנקודר מאומן. מכיוון שהמודל הגנרטיבי לתמונות הינו מודל דיפוזיה שימוש ב-VAE (חלק בלתי נפרד של מודלי דיפוזיה גנרטיביים) לא צריך להפתיע. בנוסף כאמור יש שני דקודר מאומנים שאליהם מוזנים הייצוגים הנבנים על ידי LLM.', $false, $false, $false, $false, $false, $true, 1, $false, 'השיטה פשוטה: הרכבה של פונקציית לוס משני לוסים שאחד מכם הוא ל- InTn גנרטיבי והשני ל- InTn ייצוגי. לכל אחת מהמשימות מחובר למודל ההתחלתי ראש מאומן (כמה בלוקים של טרנספורמרים למיטב הבנתי).', 2) | Out-Null
$d.Content.Find.Execute('המאמר מציע שיטה תלת שלבית לאימון המודלים כאשר כל שלב ״מפשירים״ יותר ויותר מודלים (כולל LLM) כאשר בשלב האחרון מאמנים את כולם פרט ל-VAE. ', $false, $false, $false, $false, $false, $true, 1, $false, 'אז למשימה הראשונה המחברים משתמשים בלוס הסטנדרטי של מודלי שפה גנרטיביים כלומר חיזוי של טוקן הבא עבור התשובה בלבד. למשימה השנייה המחברים משתמשים בלוס הניגודי (די סטנדרטי במשימות כאלו) והמנסה לקרב אמבדינגס של השאלה עם התשובה הנכונה ולהרחיק את האמבדינגס של השאלה עם תשובה לשאלה אחרת. ייצוג של הטקסט מחושב על ידי מודל באופן דו כיווני (אנדוקר) כאשר האמבדינג הוא הממוצע של האמבדינגס של כל הטוקנים של הטקסט. כמובן שכל משימה מקבלת פרומפט משלה.', 2) | Out-Null
$d.Content.Find.Execute('מודלי דיפוזיה במאמר מבוססי על (rectified flows (RF המנסה למפות את הדאטה מהתפלגות פשוטה (גאוסית) להתפלגות הדאטה בצורה ישרה כלומר המסלול בין x_0 הגאוסי ל x_1 של הדאטה הוא ישר. כלומר כל נקודה x_t במסלול הזה היא צירוף קמור של x_0 ו-x_1. בגדול מודל הדיפוזיה מאומן לשערך את המהירות הקבועה v(השווה x_0 - x_1 עבור כל נקודה x_t במסלול. הדגימה מבוצעת על ידי פתרון משוואה דיפרנציאלית המתארת התקדמות של x_0 ל x_1 עם מהירות v (שיטת אוילר). מודל דיפוזיה המאומן במאמר הוא לטנטי.', $false, $false, $false, $false, $false, $true, 1, $false, 'זהו זה - סקירה קלילה כמו שהבטחתי…', 2) | Out-Null

# 3. Remove the two trailing paragraphs that no longer belong in the review
#    ("פרט מעניין על המאמר..." and "מאמר כתוב יפה ודי ברור - מומלץ!")
#    They immediately follow the paragraph that used to hold the diffusion-model
#    text (now "זהו זה...") and precede the arXiv link paragraph.
$target1 = 'פרט מעניין על המאמר: אחד האיברים בפונקציית לוס של מודל דיפוזיה קונסת אותו על אי התאמה של ייצוג הפנימי המורעש (המחושב על שכבות הביניים של המודל) לייצוג התמונה הנקייה המחושב על אנקודר חזק (understanding encoder). וכמבון יש classifier guidance באימון של מודל דיפוזיה (קלאסי)'
$target2 = 'מאמר כתוב יפה ודי ברור - מומלץ!'
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $ptext = $d.Paragraphs.Item($i).Range.Text
    if ($ptext.TrimEnd() -eq $target1 -or $ptext.TrimEnd() -eq $target2) {
        $d.Paragraphs.Item($i).Range.Delete()
    }
}

# 4. Update the arxiv URL to the new paper
$d.Content.Find.Execute('https://arxiv.org/abs/2411.07975', $false, $false, $false, $false, $false, $true, 1, $false, 'https://arxiv.org/abs/2402.09906', 2) | Out-Null

Write-Host "Final paragraph count:" $d.Paragraphs.Count
